# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for each coin row per the latest scrape.
#
# Note: several Price values are strings that look like numbers (e.g.
# "1.010", "12.10") and must retain their exact text (incl. trailing
# zeros), so they're entered with a leading apostrophe to force literal
# text entry, matching the source data's text formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.710.93'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.863.88'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  +0.97%  '
$ws.Range("D5").Value = "'332.93"
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("D7").Value = "'0.4656"
$ws.Range("E7").Value = '  -1.42%  '
$ws.Range("D8").Value = "'0.3891"
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").Value = "'46.42"
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("D10").Value = "'0.07958"
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("D11").Value = "'0.9991"
$ws.Range("E11").Value = '  -2.86%  '
$ws.Range("D12").Value = "'21.50"
$ws.Range("E12").Value = '  -2.99%  '
$ws.Range("D13").Value = '1.864.68'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = "'5.986"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").Value = "'7.166"
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = "'1.012"
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = "'87.94"
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").Value = "'0.06698"
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").Value = "'16.91"
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("D22").Value = '27.690.54'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").Value = "'5.459"
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("D24").Value = "'10.87"
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").Value = "'2.322"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").Value = '2.090.57'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = "'157.91"
$ws.Range("E28").Value = '  -2.73%  '
$ws.Range("D29").Value = "'2.101"
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").Value = "'5.379"
$ws.Range("E30").Value = '  -3.98%  '
$ws.Range("D31").Value = "'121.03"
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("D32").Value = "'0.9692"
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("D33").Value = "'0.09438"
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = "'3.637"
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").Value = "'5.294"
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").Value = "'1.335"
$ws.Range("E36").Value = '  -8.29%  '
$ws.Range("D37").Value = "'0.06028"
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").Value = "'0.02208"
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("D39").Value = "'1.199"
$ws.Range("E39").Value = '  -2.80%  '
$ws.Range("D40").Value = "'8.161"
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").Value = "'1.010"
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("D42").Value = "'0.5908"
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("D43").Value = "'0.1879"
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").Value = "'10.19"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").Value = "'1.255"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").Value = "'0.5611"
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("D47").Value = "'12.10"
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("D48").Value = "'1.913"
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").Value = "'3.301"
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").Value = "'0.06762"
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("D51").Value = "'112.21"
$ws.Range("E51").Value = '  -2.14%  '
